$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (t="inlineStr" in the original file).
# Excel auto-parses a bare numeric-looking assignment (e.g. "394.87") into a
# real number, which would change the stored cell type. For the cells whose new
# price text parses as a plain number, force the cell format to Text first so
# the assignment below is kept as a literal string, matching the source data.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "57.124.77"
$ws.Range("E2").Value = "  +5.04%  "

$ws.Range("D3").Value = "3.247.49"
$ws.Range("E3").Value = "  +2.43%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "394.87"
$ws.Range("E5").Value = "  -1.02%  "

$ws.Range("D6").Value = "107.83"
$ws.Range("E6").Value = "  -1.32%  "

$ws.Range("D7").Value = "0.586"
$ws.Range("E7").Value = "  +6.72%  "

$ws.Range("D8").Value = "3.242.12"
$ws.Range("E8").Value = "  +2.48%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  +0.75%  "

$ws.Range("D11").Value = "39.06"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("E12").Value = "  +12.13%  "

$ws.Range("E13").Value = "  +1.77%  "

$ws.Range("D14").Value = "3.756.68"
$ws.Range("E14").Value = "  +2.33%  "

$ws.Range("D15").Value = "8.24"
$ws.Range("E15").Value = "  +2.60%  "

$ws.Range("D16").Value = "19.10"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").Value = "3.230.26"
$ws.Range("E17").Value = "  +1.41%  "

$ws.Range("E18").Value = "  -3.04%  "

$ws.Range("E19").Value = "  +2.40%  "

$ws.Range("D20").Value = "56.926.05"
$ws.Range("E20").Value = "  +4.74%  "

$ws.Range("D21").Value = "3.33"
$ws.Range("E21").Value = "  +0.95%  "

$ws.Range("D22").Value = "0.0000112"
$ws.Range("E22").Value = "  +13.01%  "

$ws.Range("D23").Value = "12.87"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "295.18"
$ws.Range("E24").Value = "  +8.34%  "

$ws.Range("D25").Value = "74.21"
$ws.Range("E25").Value = "  +3.70%  "

$ws.Range("D26").Value = "3.15"
$ws.Range("E26").Value = "  -2.95%  "

$ws.Range("D27").Value = "27.99"
$ws.Range("E27").Value = "  +1.19%  "

$ws.Range("D28").Value = "7.56"
$ws.Range("E28").Value = "  -5.51%  "

$ws.Range("D29").Value = "7.21"
$ws.Range("E29").Value = "  -2.05%  "

$ws.Range("D30").Value = "0.167"
$ws.Range("E30").Value = "  -1.36%  "

$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("E32").Value = "  +2.35%  "

$ws.Range("D33").Value = "0.108"
$ws.Range("E33").Value = "  -3.29%  "

$ws.Range("D34").Value = "39.47"
$ws.Range("E34").Value = "  +7.17%  "

$ws.Range("E35").Value = "  -2.91%  "

$ws.Range("E36").Value = "  +2.43%  "

$ws.Range("D37").Value = "51.41"
$ws.Range("E37").Value = "  +1.68%  "

$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("D39").Value = "3.46"
$ws.Range("E39").Value = "  -5.16%  "

$ws.Range("D40").Value = "2.90"
$ws.Range("E40").Value = "  +2.45%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "134.70"
$ws.Range("E41").Value = "  +3.51%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.122"
$ws.Range("E42").Value = "  +4.52%  "

$ws.Range("E43").Value = "  -1.92%  "

$ws.Range("D44").Value = "16.97"
$ws.Range("E44").Value = "  -1.74%  "

$ws.Range("E45").Value = "  -4.41%  "

$ws.Range("E46").Value = "  -4.34%  "

$ws.Range("D47").Value = "22.01"
$ws.Range("E47").Value = "  -0.96%  "

$ws.Range("D48").Value = "2.13"
$ws.Range("E48").Value = "  +3.10%  "

$ws.Range("D49").Value = "2.154.39"
$ws.Range("E49").Value = "  +3.19%  "

$ws.Range("E50").Value = "  -6.39%  "

$ws.Range("D51").Value = "1.93"
$ws.Range("E51").Value = "  +14.74%  "
